# Update net value for Explorer fund record
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll position / selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F20").Select()

# --- Row 41: clear stray empty cells Q41 / R41 ---
$ws.Range("Q41").ClearContents()
$ws.Range("R41").ClearContents()

# --- Row 42: clear stray empty cell I42, update O42, turn P42/S42 into formulas ---
$ws.Range("I42").ClearContents()
$ws.Range("O42").Value = 441100
$ws.Range("P42").Formula = "=Q42*R42"
$ws.Range("S42").Formula = "=T42*U42"

# --- Prepare formatting for the two new data rows (43 and 44) by copying
#     the number formats from row 42, which already has the correct layout ---
$fmtCols = @("B","C","D","E","F","G","H","M","N","O","P","S","T")
foreach ($col in $fmtCols) {
    $ws.Range($col + "42").Copy()
    $ws.Range($col + "43").PasteSpecial(-4122)
    $ws.Range($col + "44").PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = 0

# row 43 has no entry in column I, row 44 keeps an empty (unformatted) one
$ws.Range("I43").ClearFormats()
$ws.Range("I43").ClearContents()
$ws.Range("I44").ClearFormats()

# --- Row 43 data ---
$ws.Range("A43").Value = 20170707
$ws.Range("B43").Formula = "=SUM(C43:H43)"
$ws.Range("C43").Value = 1049525
$ws.Range("D43").Value = 168237
$ws.Range("E43").Value = 1215759
$ws.Range("F43").Value = 385665
$ws.Range("G43").Value = 1224413
$ws.Range("H43").Value = 1242806
$ws.Range("M43").Value = -0.4763
$ws.Range("N43").Value = 2.195
$ws.Range("O43").Value = 439000
$ws.Range("P43").Formula = "=Q43*R43"
$ws.Range("Q43").Value = 666684.18
$ws.Range("R43").Value = 2.9824
$ws.Range("S43").Formula = "=T43*U43"
$ws.Range("T43").Value = 1748728.62
$ws.Range("U43").Value = 1.6349

# --- Row 44 data ---
$ws.Range("A44").Value = 20170714
$ws.Range("B44").Formula = "=SUM(C44:H44)"
$ws.Range("C44").Value = 1018647
$ws.Range("D44").Value = 168237
$ws.Range("E44").Value = 1222601
$ws.Range("F44").Value = 385429
$ws.Range("G44").Value = 1248767
$ws.Range("H44").Value = 1255327
$ws.Range("M44").Value = 0.2384
$ws.Range("N44").Value = 2.2
$ws.Range("O44").Value = 440000
$ws.Range("P44").Formula = "=Q44*R44"
$ws.Range("Q44").Value = 666684.18
$ws.Range("R44").Value = 2.9895
$ws.Range("S44").Formula = "=T44*U44"
$ws.Range("T44").Value = 1748728.62
$ws.Range("U44").Value = 1.6388
